# "Add files via upload" — re-upload of the Test-Cases tracker with the
# "Approved/Rejected" column (I) filled in with "Approved" for every test
# step row (rows 2-23), plus the sheet's scroll position / active
# selection moved to reflect where the author was last looking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 9).Value = "Approved"
}

# Scroll the view so column H is the leftmost visible column (topLeftCell
# H1) and select J10, matching the saved window state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("J10").Select()
